# Apply weekly update: shift existing Hortaliza/Sandia data rows down by 2
# and insert two new rows at the top of the data block (rows 26-27),
# extending the data block from row 108 to row 110.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, D(Fecha), I(Calidad), J(Volumen), K(Precio minimo), L(Precio maximo),
#             M(Precio promedio ponderado), N(Unidad de comercializacion), O(Origen), P(Precio $/Kg)
$rows = @(
    @(26, 44550, 'Extra', 3000, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(27, 44550, 'Primera', 2000, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(28, 44547, 'Extra', 3000, 3000, 3000, 3000, '$/unidad', 'Región de O''Higgins', 3000),
    @(29, 44547, 'Extra', 3000, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(30, 44547, 'Primera', 3000, 2500, 2500, 2500, '$/unidad', 'Región de O''Higgins', 2500),
    @(31, 44547, 'Primera', 3000, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(32, 44223, 'Extra', 2000, 2300, 2300, 2300, '$/unidad', 'Región del Maule', 2300),
    @(33, 44223, 'Primera', 5000, 1800, 1800, 1800, '$/unidad', 'Región del Maule', 1800),
    @(34, 44223, 'Segunda', 3000, 1300, 1300, 1300, '$/unidad', 'Región del Maule', 1300),
    @(35, 44202, 'Extra', 5000, 2300, 2300, 2300, '$/unidad', 'Región del Maule', 2300),
    @(36, 44202, 'Primera', 10000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(37, 44202, 'Segunda', 8000, 1800, 1800, 1800, '$/unidad', 'Región del Maule', 1800),
    @(38, 44216, 'Extra', 3000, 1800, 1800, 1800, '$/unidad', 'Región del Maule', 1800),
    @(39, 44216, 'Primera', 5000, 1400, 1400, 1400, '$/unidad', 'Región del Maule', 1400),
    @(40, 44216, 'Segunda', 3000, 1000, 1000, 1000, '$/unidad', 'Región del Maule', 1000),
    @(41, 44186, 'Extra', 2500, 3200, 3200, 3200, '$/unidad', 'Región del Maule', 3200),
    @(42, 44186, 'Primera', 3500, 2800, 2800, 2800, '$/unidad', 'Región del Maule', 2800),
    @(43, 44189, 'Extra', 1000, 3500, 3500, 3500, '$/unidad', 'Región del Maule', 3500),
    @(44, 44189, 'Primera', 2500, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(45, 44189, 'Segunda', 1500, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(46, 44169, 'Primera', 2000, 400, 400, 400, '$/kilo (volumen en unidades)', 'Región del Maule', 400),
    @(47, 44215, 'Extra', 5000, 1800, 1800, 1800, '$/unidad', 'Región del Maule', 1800),
    @(48, 44215, 'Primera', 9000, 1500, 1500, 1500, '$/unidad', 'Región del Maule', 1500),
    @(49, 44215, 'Segunda', 6000, 1200, 1200, 1200, '$/unidad', 'Región del Maule', 1200),
    @(50, 44201, 'Extra', 5000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(51, 44201, 'Primera', 10000, 1700, 1700, 1700, '$/unidad', 'Región del Maule', 1700),
    @(52, 44201, 'Segunda', 7000, 1200, 1200, 1200, '$/unidad', 'Región del Maule', 1200),
    @(53, 44179, 'Especial', 2000, 3200, 3200, 3200, '$/unidad', 'Región del Maule', 3200),
    @(54, 44179, 'Primera', 3000, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(55, 44203, 'Extra', 3000, 2200, 2200, 2200, '$/unidad', 'Región del Maule', 2200),
    @(56, 44203, 'Primera', 5000, 1800, 1800, 1800, '$/unidad', 'Región del Maule', 1800),
    @(57, 44203, 'Segunda', 4000, 1300, 1300, 1300, '$/unidad', 'Región del Maule', 1300),
    @(58, 44214, 'Extra', 3000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(59, 44214, 'Primera', 8000, 1600, 1600, 1600, '$/unidad', 'Región del Maule', 1600),
    @(60, 44214, 'Segunda', 4000, 1200, 1200, 1200, '$/unidad', 'Región del Maule', 1200),
    @(61, 44546, 'Extra', 2000, 3000, 3000, 3000, '$/unidad', 'Región de O''Higgins', 3000),
    @(62, 44546, 'Extra', 2500, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(63, 44546, 'Primera', 3000, 2500, 2500, 2500, '$/unidad', 'Región de O''Higgins', 2500),
    @(64, 44546, 'Primera', 3500, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(65, 44222, 'Extra', 3000, 2300, 2300, 2300, '$/unidad', 'Región del Maule', 2300),
    @(66, 44222, 'Primera', 6000, 1800, 1800, 1800, '$/unidad', 'Región del Maule', 1800),
    @(67, 44222, 'Segunda', 4000, 1300, 1300, 1300, '$/unidad', 'Región del Maule', 1300),
    @(68, 44181, 'Primera', 5000, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(69, 44181, 'Segunda', 3000, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(70, 44209, 'Extra', 3000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(71, 44209, 'Primera', 5000, 1600, 1600, 1600, '$/unidad', 'Región del Maule', 1600),
    @(72, 44209, 'Segunda', 5000, 1200, 1200, 1200, '$/unidad', 'Región del Maule', 1200),
    @(73, 44533, 'Primera', 2000, 500, 500, 500, '$/kilo', 'Perú', 500),
    @(74, 44176, 'Primera', 1500, 400, 400, 400, '$/kilo (volumen en unidades)', 'Perú', 400),
    @(75, 44176, 'Primera', 1500, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(76, 44176, 'Segunda', 1500, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(77, 44210, 'Extra', 3000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(78, 44210, 'Primera', 6000, 1600, 1600, 1600, '$/unidad', 'Región del Maule', 1600),
    @(79, 44210, 'Segunda', 4000, 1200, 1200, 1200, '$/unidad', 'Región del Maule', 1200),
    @(80, 44168, 'Primera', 1200, 450, 450, 450, '$/kilo (volumen en unidades)', 'Región Metropolitana', 450),
    @(81, 44231, 'Extra', 3000, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(82, 44231, 'Primera', 6000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(83, 44231, 'Segunda', 4000, 1500, 1500, 1500, '$/unidad', 'Región del Maule', 1500),
    @(84, 44208, 'Extra', 4000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(85, 44208, 'Primera', 5000, 1800, 1800, 1800, '$/unidad', 'Región del Maule', 1800),
    @(86, 44208, 'Segunda', 3000, 1200, 1200, 1200, '$/unidad', 'Región del Maule', 1200),
    @(87, 44264, 'Extra', 2000, 2300, 2300, 2300, '$/unidad', 'Región del Maule', 2300),
    @(88, 44264, 'Primera', 2000, 1800, 1800, 1800, '$/unidad', 'Región del Maule', 1800),
    @(89, 44264, 'Segunda', 2000, 1300, 1300, 1300, '$/unidad', 'Región del Maule', 1300),
    @(90, 44232, 'Extra', 2000, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(91, 44232, 'Primera', 5000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(92, 44232, 'Segunda', 3000, 1600, 1600, 1600, '$/unidad', 'Región del Maule', 1600),
    @(93, 44196, 'Extra', 1500, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(94, 44196, 'Primera', 2500, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(95, 44196, 'Segunda', 2000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(96, 44200, 'Extra', 4000, 2200, 2200, 2200, '$/unidad', 'Región del Maule', 2200),
    @(97, 44200, 'Primera', 8000, 1700, 1700, 1700, '$/unidad', 'Región del Maule', 1700),
    @(98, 44200, 'Segunda', 4000, 1400, 1400, 1400, '$/unidad', 'Región del Maule', 1400),
    @(99, 44188, 'Extra', 3000, 3500, 3500, 3500, '$/unidad', 'Región del Maule', 3500),
    @(100, 44188, 'Primera', 3000, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(101, 44188, 'Segunda', 3000, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(102, 44224, 'Extra', 2000, 2300, 2300, 2300, '$/unidad', 'Región del Maule', 2300),
    @(103, 44224, 'Primera', 6000, 1800, 1800, 1800, '$/unidad', 'Región del Maule', 1800),
    @(104, 44224, 'Segunda', 4000, 1300, 1300, 1300, '$/unidad', 'Región del Maule', 1300),
    @(105, 44195, 'Extra', 2000, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(106, 44195, 'Primera', 3000, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(107, 44195, 'Segunda', 2000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000),
    @(108, 44194, 'Extra', 2000, 3000, 3000, 3000, '$/unidad', 'Región del Maule', 3000),
    @(109, 44194, 'Primera', 4000, 2500, 2500, 2500, '$/unidad', 'Región del Maule', 2500),
    @(110, 44194, 'Segunda', 2000, 2000, 2000, 2000, '$/unidad', 'Región del Maule', 2000)
)

foreach ($row in $rows) {
    $r = $row[0]

    # Columns constant across this data block (unchanged by the edit)
    $ws.Cells.Item($r, 1).Value = 5
    $ws.Cells.Item($r, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($r, 3).Value = "Maule"

    # D: Fecha
    $ws.Cells.Item($r, 4).Value = $row[1]

    $ws.Cells.Item($r, 5).Value = 7
    $ws.Cells.Item($r, 6).Value = 100112028
    $ws.Cells.Item($r, 7).Value = "Sandia"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"

    # I: Calidad
    $ws.Cells.Item($r, 9).Value = $row[2]
    # J: Volumen
    $ws.Cells.Item($r, 10).Value = $row[3]
    # K: Precio minimo
    $ws.Cells.Item($r, 11).Value = $row[4]
    # L: Precio maximo
    $ws.Cells.Item($r, 12).Value = $row[5]
    # M: Precio promedio ponderado
    $ws.Cells.Item($r, 13).Value = $row[6]
    # N: Unidad de comercializacion
    $ws.Cells.Item($r, 14).Value = $row[7]
    # O: Origen
    $ws.Cells.Item($r, 15).Value = $row[8]
    # P: Precio $/Kg
    $ws.Cells.Item($r, 16).Value = $row[9]

    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
